$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - copy G1's formatting (bold, centered, bordered
# header style) to H1, then overwrite the copied value with the new header text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column in row 2
$ws.Range("H2").Value = 0
